{"js": "// Apply the 26 ordered text replacements (1 date header + 25 multiplication\n// problems) described by the diff. Every `old` value is unique within the\n// document, so a straightforward search-and-replace per pair is safe and\n// order-independent.\nconst replacements = [\n  [\"2025-05-18 Sunday\", \"2025-05-19 Monday\"],\n  [\"461\u00d72=\", \"781\u00d78=\"],\n  [\"118\u00d75=\", \"542\u00d78=\"],\n  [\"989\u00d73=\", \"695\u00d72=\"],\n  [\"919\u00d74=\", \"167\u00d76=\"],\n  [\"951\u00d77=\", \"438\u00d77=\"],\n  [\"728\u00d73=\", \"328\u00d75=\"],\n  [\"163\u00d73=\", \"204\u00d77=\"],\n  [\"512\u00d79=\", \"793\u00d74=\"],\n  [\"483\u00d73=\", \"267\u00d78=\"],\n  [\"762\u00d72=\", \"370\u00d76=\"],\n  [\"731\u00d78=\", \"641\u00d72=\"],\n  [\"748\u00d77=\", \"168\u00d74=\"],\n  [\"427\u00d72=\", \"497\u00d79=\"],\n  [\"224\u00d74=\", \"408\u00d78=\"],\n  [\"436\u00d76=\", \"377\u00d78=\"],\n  [\"988\u00d78=\", \"583\u00d79=\"],\n  [\"609\u00d73=\", \"834\u00d77=\"],\n  [\"138\u00d78=\", \"778\u00d75=\"],\n  [\"334\u00d77=\", \"624\u00d75=\"],\n  [\"622\u00d77=\", \"176\u00d76=\"],\n  [\"334\u00d72=\", \"188\u00d76=\"],\n  [\"791\u00d75=\", \"438\u00d79=\"],\n  [\"993\u00d78=\", \"317\u00d77=\"],\n  [\"895\u00d75=\", \"585\u00d78=\"],\n  [\"856\u00d75=\", \"643\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the 26 ordered text replacements (1 date header + 25 multiplication\n# problems) described by the diff. Every `old` value is unique within the\n# document, so Find/Replace per pair (restricted to exact match, no wildcards)\n# is safe and order-independent.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2025-05-18 Sunday', '2025-05-19 Monday'),\n    @('461\u00d72=', '781\u00d78='),\n    @('118\u00d75=', '542\u00d78='),\n    @('989\u00d73=', '695\u00d72='),\n    @('919\u00d74=', '167\u00d76='),\n    @('951\u00d77=', '438\u00d77='),\n    @('728\u00d73=', '328\u00d75='),\n    @('163\u00d73=', '204\u00d77='),\n    @('512\u00d79=', '793\u00d74='),\n    @('483\u00d73=', '267\u00d78='),\n    @('762\u00d72=', '370\u00d76='),\n    @('731\u00d78=', '641\u00d72='),\n    @('748\u00d77=', '168\u00d74='),\n    @('427\u00d72=', '497\u00d79='),\n    @('224\u00d74=', '408\u00d78='),\n    @('436\u00d76=', '377\u00d78='),\n    @('988\u00d78=', '583\u00d79='),\n    @('609\u00d73=', '834\u00d77='),\n    @('138\u00d78=', '778\u00d75='),\n    @('334\u00d77=', '624\u00d75='),\n    @('622\u00d77=', '176\u00d76='),\n    @('334\u00d72=', '188\u00d76='),\n    @('791\u00d75=', '438\u00d79='),\n    @('993\u00d78=', '317\u00d77='),\n    @('895\u00d75=', '585\u00d78='),\n    @('856\u00d75=', '643\u00d77='),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
